$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 — append the 2026-02-19 entry (class cancelled for a "reuniao"
# meeting, per the calendar-rules integration described in the commit).
# Force the whole row to Text first so values that look like dates/numbers
# (the date itself, "26.0", "15:15") are stored as literal strings, matching
# how every other row in this sheet already stores its data (inlineStr/text).
$row = $ws.Range("A18:M18")
$row.NumberFormat = "@"

$ws.Range("A18").Value = "2026-02-19"
$ws.Range("B18").Value = "Nublado"
$ws.Range("C18").Value = "Agradavel"
$ws.Range("D18").Value = "cancelada"
$ws.Range("E18").Value = "reuniao"
$ws.Range("F18").Value = "nenhuma"
$ws.Range("G18").Value = "26.0"
$ws.Range("H18").Value = "'"
$ws.Range("I18").Value = "-"
$ws.Range("J18").Value = "jtq07"
$ws.Range("K18").Value = "Terça e Quinta"
$ws.Range("L18").Value = "15:15"
$ws.Range("M18").Value = "Jefferson"

# Drop the explicit text format again so the new cells end up unstyled, same
# as the rest of the sheet (no "s" attribute / style reference).
$row.ClearFormats()
